$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.492.62"
$ws.Range("E2").Value = "  +2.94%  "
$ws.Range("D3").Value = "1.603.39"
$ws.Range("E3").Value = "  +2.46%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.519"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "26.53"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.34"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.250"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0598"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.25%  "
$ws.Range("E12").Value = "  +1.68%  "
$ws.Range("D13").Value = "1.834.45"
$ws.Range("E13").Value = "  +2.58%  "
$ws.Range("D14").Value = "1.602.03"
$ws.Range("E14").Value = "  +2.38%  "
$ws.Range("D15").Value = "29.518.12"
$ws.Range("E15").Value = "  +2.98%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.532"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.57%  "
$ws.Range("D21").Value = "0.0₃0689"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("E23").Value = "  +1.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.27%  "
$ws.Range("E28").Value = "  +5.02%  "
$ws.Range("E29").Value = "  +2.28%  "
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("E31").Value = "  +2.44%  "
$ws.Range("E32").Value = "  -0.37%  "
$ws.Range("E33").Value = "  +1.62%  "
$ws.Range("E34").Value = "  +3.83%  "
$ws.Range("D35").Value = "1.415.08"
$ws.Range("E35").Value = "  +1.79%  "
$ws.Range("E36").Value = "  -1.09%  "
$ws.Range("E37").Value = "  +2.86%  "
$ws.Range("E38").Value = "  +5.46%  "
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("E40").Value = "  +1.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.535"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.40%  "
$ws.Range("E42").Value = "  +0.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "53.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +24.38%  "
$ws.Range("E44").Value = "  +5.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.998"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.793"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.13%  "
$ws.Range("E47").Value = "  +2.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "1.745.10"
$ws.Range("E49").Value = "  +2.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "86.65"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.843"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.95%  "
